$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5342
$ws1.Range("F5").Value = 7573
$ws1.Range("F6").Value = 48
$ws1.Range("F9").Value = 606
$ws1.Range("F11").Value = 34
$ws1.Range("F12").Value = 4373
$ws1.Range("F13").Value = 1781
$ws1.Range("F16").Value = 2961
$ws1.Range("F18").Value = 570
$ws1.Range("F19").Value = 216
$ws1.Range("F20").Value = 534
$ws1.Range("F21").Value = 463
$ws1.Range("F22").Value = 469
$ws1.Range("F23").Value = 328
$ws1.Range("F24").Value = 111
$ws1.Range("F25").Value = 1713
$ws1.Range("F26").Value = 1222
$ws1.Range("F27").Value = 98
$ws1.Range("F28").Value = 1410
$ws1.Range("F29").Value = 116
$ws1.Range("F30").Value = 590
$ws1.Range("F32").Value = 518
$ws1.Range("F37").Value = 72
$ws1.Range("F38").Value = 3021
$ws1.Range("F41").Value = 123
$ws1.Range("F42").Value = 45
$ws1.Range("F43").Value = 252

# Sheet "全部类型" (All types) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5342
$ws4.Range("F5").Value = 7573
$ws4.Range("F6").Value = 48
$ws4.Range("F9").Value = 606
$ws4.Range("F11").Value = 34
$ws4.Range("F12").Value = 4373
$ws4.Range("F13").Value = 1781
$ws4.Range("F16").Value = 2961
$ws4.Range("F18").Value = 570
$ws4.Range("F19").Value = 216
$ws4.Range("F20").Value = 534
$ws4.Range("F21").Value = 463
$ws4.Range("F22").Value = 469
$ws4.Range("F24").Value = 328
$ws4.Range("F25").Value = 111
$ws4.Range("F26").Value = 1713
$ws4.Range("F27").Value = 1222
$ws4.Range("F28").Value = 98
$ws4.Range("F29").Value = 1410
$ws4.Range("F30").Value = 116
$ws4.Range("F31").Value = 590
$ws4.Range("F33").Value = 518
$ws4.Range("F38").Value = 72
$ws4.Range("F39").Value = 3021
$ws4.Range("F43").Value = 123
$ws4.Range("F44").Value = 45
$ws4.Range("F45").Value = 253
